$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.854.89"
$ws.Range("E2").Value = "'  +6.93%  "
$ws.Range("D3").Value = "'3.019.00"
$ws.Range("E3").Value = "'  +4.31%  "
$ws.Range("E4").Value = "'  +0.17%  "
$ws.Range("D5").Value = "'585.37"
$ws.Range("E5").Value = "'  +3.07%  "
$ws.Range("D6").Value = "'154.79"
$ws.Range("E6").Value = "'  +7.75%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "'  -0.20%  "
$ws.Range("D8").Value = "'3.014.16"
$ws.Range("E8").Value = "'  +4.26%  "
$ws.Range("D9").Value = "'0.519"
$ws.Range("E9").Value = "'  +2.64%  "
$ws.Range("D10").Value = "'6.98"
$ws.Range("E10").Value = "'  -0.20%  "
$ws.Range("E11").Value = "'  +5.00%  "
$ws.Range("D12").Value = "'0.452"
$ws.Range("E12").Value = "'  +5.03%  "
$ws.Range("E13").Value = "'  +3.74%  "
$ws.Range("D14").Value = "'34.34"
$ws.Range("E14").Value = "'  +7.66%  "
$ws.Range("D15").Value = "'0.127"
$ws.Range("E15").Value = "'  +0.93%  "
$ws.Range("D16").Value = "'65.889.15"
$ws.Range("E16").Value = "'  +6.98%  "
$ws.Range("D17").Value = "'3.517.98"
$ws.Range("E17").Value = "'  +4.24%  "
$ws.Range("D18").Value = "'6.99"
$ws.Range("E18").Value = "'  +6.61%  "
$ws.Range("D19").Value = "'3.019.52"
$ws.Range("E19").Value = "'  +3.95%  "
$ws.Range("D20").Value = "'460.60"
$ws.Range("E20").Value = "'  +6.43%  "
$ws.Range("E21").Value = "'  +6.03%  "
$ws.Range("D22").Value = "'0.687"
$ws.Range("E22").Value = "'  +4.80%  "
$ws.Range("D23").Value = "'7.40"
$ws.Range("E23").Value = "'  +8.58%  "
$ws.Range("D24").Value = "'82.03"
$ws.Range("E24").Value = "'  +3.46%  "
$ws.Range("D25").Value = "'12.67"
$ws.Range("E25").Value = "'  +5.62%  "
$ws.Range("E26").Value = "'  +12.07%  "
$ws.Range("D27").Value = "'10.79"
$ws.Range("E27").Value = "'  +8.13%  "
$ws.Range("E28").Value = "'  +0.00%  "
$ws.Range("D29").Value = "'2.42"
$ws.Range("E29").Value = "'  +18.37%  "
$ws.Range("D30").Value = "'7.92"
$ws.Range("E30").Value = "'  +12.89%  "
$ws.Range("E31").Value = "'  +4.25%  "
$ws.Range("E32").Value = "'  -2.83%  "
$ws.Range("D33").Value = "'0.113"
$ws.Range("E33").Value = "'  +5.98%  "
$ws.Range("D34").Value = "'27.07"
$ws.Range("E34").Value = "'  +6.22%  "
$ws.Range("E35").Value = "'  +0.07%  "
$ws.Range("D36").Value = "'0.992"
$ws.Range("E36").Value = "'  +3.59%  "
$ws.Range("D37").Value = "'5.82"
$ws.Range("E37").Value = "'  +7.93%  "
$ws.Range("E38").Value = "'  +11.27%  "
$ws.Range("B39").Value = "'dogwifhat"
$ws.Range("C39").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "'2.98"
$ws.Range("E39").Value = "'  +6.00%  "
$ws.Range("B40").Value = "'Arweave"
$ws.Range("C40").Value = "'https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D40").Value = "'45.39"
$ws.Range("E40").Value = "'  +15.12%  "
$ws.Range("B41").Value = "'OKB"
$ws.Range("C41").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Value = "'49.48"
$ws.Range("E41").Value = "'  +1.34%  "
$ws.Range("E42").Value = "'  +6.70%  "
$ws.Range("D43").Value = "'0.303"
$ws.Range("E43").Value = "'  +13.76%  "
$ws.Range("D44").Value = "'8.48"
$ws.Range("E44").Value = "'  +3.35%  "
$ws.Range("D45").Value = "'387.98"
$ws.Range("E45").Value = "'  +11.96%  "
$ws.Range("D46").Value = "'2.794.02"
$ws.Range("E46").Value = "'  +3.54%  "
$ws.Range("D47").Value = "'0.0354"
$ws.Range("E47").Value = "'  +5.74%  "
$ws.Range("D48").Value = "'134.35"
$ws.Range("E48").Value = "'  +1.39%  "
$ws.Range("D50").Value = "'23.64"
$ws.Range("E50").Value = "'  +9.51%  "
$ws.Range("D51").Value = "'0.107"
$ws.Range("E51").Value = "'  +3.79%  "
